$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (Cylinder / Engine / N / Mean / SD / Mean / SD) live in A5:G7
# (row 8 - the 8-cylinder group - already sorts after them and is untouched).
# Sort those three rows ascending by column A ("Cylinder") so 4 < 6 < 6.
$dataRange = $ws.Range("A5:G7")
$sortKey = $ws.Range("A5:A7")
$dataRange.Sort($sortKey, 1)

# After sorting, rows 6 and 7 both have "6" in column A. Collapse that into a
# single merged cell instead of repeating the value, with the label aligned
# to the top of the merged block.
$ws.Range("A6:A7").Merge()

# The now-empty lower half (A7) becomes a plain, unfilled spacer cell
# (matching the look of column H), so copy that formatting over and clear
# the leftover value.
$ws.Range("H5").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = ""

$ws.Range("A6").VerticalAlignment = -4160
